$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test-case rows (39-43): WithDraw Deal Testcase Update ---
# Values are written row-by-row, column A then column D for each row
# (across all five rows), and finally column B for all five rows, so new
# shared-string entries land in the same order the original authoring
# tool produced them in.

# Columns A (Automation Test ID) and D (Expected Result), row by row
$ws.Range("A39").Value = "DuplicateDeal_TC001"
$ws.Range("D39").Value = "Deal has been duplicated"
$ws.Range("A40").Value = "DuplicateDeal_TC002"
$ws.Range("D40").Value = "Deal has been duplicated"
$ws.Range("A41").Value = "WithDrawDeal_TC001"
$ws.Range("D41").Value = "Deal has been withdrawn"
$ws.Range("A42").Value = "WithDrawDeal_TC002"
$ws.Range("D42").Value = "Deal has been withdrawn"
$ws.Range("A43").Value = "WithDrawDeal_TC003"
$ws.Range("D43").Value = "Deal has been withdrawn"

# Column B - Description
$ws.Range("B39").Value = "Validate whether a carrier user can duplicate deal in Opportunity.`n1) Enter valid user id and Password and click Login.`n2)Create deal and share with carrier user.`n3) Login as carrier user.`n4) Check shared deal in opportunity.`n5) Duplicate deal and check whehter another deal created.`n"
$ws.Range("B40").Value = "Validate whether a Shipper user can duplicate deal in Draft.`n1) Enter valid user id and Password and click Login.`n2)Click on Deals menu.`n3)Select any deal in Opportunity/Drafts/Booked/Withdrawn tab &  click ellipsis icon and click Duplicate option.`n4) Now Duplicate popup message has been displayed.`n5) Enter number of deals to copied in pop up window.`n6) Click on Submit button .`n"
$ws.Range("B41").Value = "Validate whether a Shipper user can withdraw deal in Opportunity.`n1) Enter valid user id and Password and click Login(Shipper user).`n2)Create deal and share with carrier user.`n3) Check shared deal in opportunity.`n4) withdraw deal and check whehter deal is displayed in withdraw tab.`n"
$ws.Range("B42").Value = "Validate whether a Shipper admin user can withdraw deal in Opportunity.`n1) Enter valid user id and Password and click Login(Shipper admin).`n2)Create deal and share with all user.`n3) Check shared deal in opportunity.`n4) withdraw deal and check whehter deal is displayed in withdraw tab.`n"
$ws.Range("B43").Value = "Validate whether a carrier user can withdraw deal in Opportunity.`n1) Enter valid user id and Password and click Login.`n2)Create deal and share with all user.`n3) Check shared deal in opportunity.`n4) withdraw deal and check whehter deal is displayed in withdraw tab.`n"

# Column C - Run Mode
$ws.Range("C39").Value = "NO"
$ws.Range("C40").Value = "NO"
$ws.Range("C41").Value = "NO"
$ws.Range("C42").Value = "NO"
$ws.Range("C43").Value = "Yes"

# Wrap/valign formatting matching the rest of the sheet's data rows.
$ws.Range("A39:A43").VerticalAlignment = -4108
$ws.Range("C39:C43").VerticalAlignment = -4108
$ws.Range("D39:D43").VerticalAlignment = -4108
$ws.Range("B39:B43").WrapText = $true

# Row heights to fit the wrapped description text (matches target sheet).
$ws.Rows.Item(39).RowHeight = 120
$ws.Rows.Item(40).RowHeight = 135
$ws.Rows.Item(41).RowHeight = 135
$ws.Rows.Item(42).RowHeight = 135
$ws.Rows.Item(43).RowHeight = 120

# Update view: scroll/select near the newly-added rows.
$ws.Range("B45").Select()
$ws.Application.ActiveWindow.Zoom = 100

Write-Host "WithDraw Deal testcases added"
